# Fruta / hortaliza, semanal
# Weekly refresh of the Guayaba subset sheet: a new sampling date (44798)
# was inserted into the series (as "Primera"/"Segunda" rows), pushing the
# rest of the data table down by one row and growing the used range from
# A1:T49 to A1:T50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 50: new row, cloned from the (old) last row (row 49) template,
#     since every row shares the same static columns A,B,C,E,F,G,H,I,J,K,Q,R,T.
$ws.Range("A50").Value = $ws.Range("A49").Value2
$ws.Range("B50").Value = $ws.Range("B49").Value2
$ws.Range("C50").Value = $ws.Range("C49").Value2
$ws.Range("D50").Value = $ws.Range("D49").Value2
$ws.Range("D50").NumberFormat = $ws.Range("D49").NumberFormat
$ws.Range("E50").Value = $ws.Range("E49").Value2
$ws.Range("F50").Value = $ws.Range("F49").Value2
$ws.Range("G50").Value = $ws.Range("G49").Value2
$ws.Range("H50").Value = $ws.Range("H49").Value2
$ws.Range("I50").Value = $ws.Range("I49").Value2
$ws.Range("J50").Value = $ws.Range("J49").Value2
$ws.Range("K50").Value = $ws.Range("K49").Value2
$ws.Range("Q50").Value = $ws.Range("Q49").Value2
$ws.Range("R50").Value = $ws.Range("R49").Value2
$ws.Range("T50").Value = $ws.Range("T49").Value2

# Now fill in the target per-row values (D date, L quality, M volume,
# N/O/P/S prices) for rows 42-50, reflecting the reshuffled weekly data.

# Row 42
$ws.Range("D42").Value = 44778
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 700
$ws.Range("O42").Value = 800
$ws.Range("P42").Value = 750
$ws.Range("S42").Value = 750

# Row 43
$ws.Range("D43").Value = 44778
$ws.Range("L43").Value = "Segunda"
$ws.Range("M43").Value = 140
$ws.Range("N43").Value = 500
$ws.Range("O43").Value = 600
$ws.Range("P43").Value = 550
$ws.Range("S43").Value = 550

# Row 44
$ws.Range("D44").Value = 44798
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 100
$ws.Range("N44").Value = 700
$ws.Range("O44").Value = 800
$ws.Range("P44").Value = 750
$ws.Range("S44").Value = 750

# Row 45
$ws.Range("D45").Value = 44798
$ws.Range("L45").Value = "Segunda"
$ws.Range("M45").Value = 130
$ws.Range("N45").Value = 500
$ws.Range("O45").Value = 600
$ws.Range("P45").Value = 550
$ws.Range("S45").Value = 550

# Row 46
$ws.Range("D46").Value = 44309
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 160
$ws.Range("N46").Value = 1400
$ws.Range("O46").Value = 1500
$ws.Range("P46").Value = 1450
$ws.Range("S46").Value = 1450

# Row 47
$ws.Range("D47").Value = 44379
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 150
$ws.Range("N47").Value = 700
$ws.Range("O47").Value = 800
$ws.Range("P47").Value = 747
$ws.Range("S47").Value = 747

# Row 48
$ws.Range("D48").Value = 44379
$ws.Range("L48").Value = "Segunda"
$ws.Range("M48").Value = 140
$ws.Range("N48").Value = 500
$ws.Range("O48").Value = 600
$ws.Range("P48").Value = 543
$ws.Range("S48").Value = 543

# Row 49
$ws.Range("D49").Value = 44344
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 140
$ws.Range("N49").Value = 1000
$ws.Range("O49").Value = 1200
$ws.Range("P49").Value = 1100
$ws.Range("S49").Value = 1100

# Row 50
$ws.Range("D50").Value = 44344
$ws.Range("L50").Value = "Segunda"
$ws.Range("M50").Value = 120
$ws.Range("N50").Value = 800
$ws.Range("O50").Value = 850
$ws.Range("P50").Value = 825
$ws.Range("S50").Value = 825
